$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 5 (pushes existing rows 5-21 down to 7-23)
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# --- New row 5 ---
$ws.Cells.Item(5,1).Value2  = 5
$ws.Cells.Item(5,2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(5,3).Value2  = "Maule"
$ws.Cells.Item(5,4).Value2  = 44462
$ws.Cells.Item(5,5).Value2  = 7
$ws.Cells.Item(5,6).Value2  = "Fruta"
$ws.Cells.Item(5,7).Value2  = 100107
$ws.Cells.Item(5,8).Value2  = "Otros"
$ws.Cells.Item(5,9).Value2  = 100107002
$ws.Cells.Item(5,10).Value2 = "Chirimoya"
$ws.Cells.Item(5,11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(5,12).Value2 = "Especial"
$ws.Cells.Item(5,13).Value2 = 205
$ws.Cells.Item(5,14).Value2 = 30000
$ws.Cells.Item(5,15).Value2 = 30000
$ws.Cells.Item(5,16).Value2 = 30000
$ws.Cells.Item(5,17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(5,18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(5,19).Value2 = 3000
$ws.Cells.Item(5,20).Value2 = 10

# --- New row 6 ---
$ws.Cells.Item(6,1).Value2  = 5
$ws.Cells.Item(6,2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(6,3).Value2  = "Maule"
$ws.Cells.Item(6,4).Value2  = 44462
$ws.Cells.Item(6,5).Value2  = 7
$ws.Cells.Item(6,6).Value2  = "Fruta"
$ws.Cells.Item(6,7).Value2  = 100107
$ws.Cells.Item(6,8).Value2  = "Otros"
$ws.Cells.Item(6,9).Value2  = 100107002
$ws.Cells.Item(6,10).Value2 = "Chirimoya"
$ws.Cells.Item(6,11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(6,12).Value2 = "Primera"
$ws.Cells.Item(6,13).Value2 = 180
$ws.Cells.Item(6,14).Value2 = 28000
$ws.Cells.Item(6,15).Value2 = 28000
$ws.Cells.Item(6,16).Value2 = 28000
$ws.Cells.Item(6,17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(6,18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(6,19).Value2 = 2800
$ws.Cells.Item(6,20).Value2 = 10
